# Bài 18 : Excel VSTO - API
# Insert a new title row above the existing header row, add a new
# "Hình Ảnh" (Image) header column, and apply some formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing header row (Ngay/Gio/...) moves
# down from row 1 to row 2.
$ws.Rows.Item(1).Insert()

# New title cell in the freshly inserted row 1.
$ws.Range("A1").Value = "Thành phố"
$ws.Range("A1").Interior.Color = 65535
$ws.Range("A1").Font.Bold = $true

$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Color = 255
$ws.Range("B1").Borders.LineStyle = 1

# New header column for images, appended after "Thời tiết".
$ws.Range("F2").Value = "Hình Ảnh"

# Row heights / column formats akin to the target workbook.
$ws.Rows.Item(1).RowHeight = 26.25
$ws.Rows.Item(2).RowHeight = 26.25

$ws.Columns.Item(1).NumberFormat = "m/d/yyyy"
$ws.Columns.Item(2).NumberFormat = "[$-F400]h:mm:ss AM/PM"

$ws.Range("B2").Select()
